$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ------------------------------------------------------------------
# 1) Insert a new row for "CONTROLOC" right before the current DAFLON
#    row (row 9), so DAFLON (and everything after it) shifts down by 1.
# ------------------------------------------------------------------
$ws.Rows("9").Insert()
# Copy the formatting from the row that just got pushed down (the old
# row 9 / DAFLON, now at row 10) so the new row looks like the others.
$ws.Range("A10:Q10").Copy()
$ws.Range("A9:Q9").PasteSpecial($xlFormats)
$ws.Rows("9").RowHeight = $ws.Rows("10").RowHeight
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "CONTROLOC 20MG 14  GASTRO RESISTANT TABS"
$ws.Range("H9").Value = "0:0"

$fmt = $ws.Range("L9").NumberFormat
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "1"
$ws.Range("L9").NumberFormat = $fmt

$ws.Range("N9").Value = "188.00"

$fmt = $ws.Range("P9").NumberFormat
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "94.0000"
$ws.Range("P9").NumberFormat = $fmt

$ws.Range("Q9").Value = "0:1"

# Renumber the rows that followed (DAFLON .. سرنجات), which are now one
# row further down (rows 10-13 instead of 9-12).
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6

# ------------------------------------------------------------------
# 2) Insert a new row for "MEBEFAC" right before the PRONTOGEST row,
#    which after the step above now sits at row 12.
# ------------------------------------------------------------------
$ws.Rows("12").Insert()
$ws.Range("A13:Q13").Copy()
$ws.Range("A12:Q12").PasteSpecial($xlFormats)
$ws.Rows("12").RowHeight = $ws.Rows("13").RowHeight
$ws.Range("A12:B12").Merge()
$ws.Range("C12:G12").Merge()
$ws.Range("H12:K12").Merge()
$ws.Range("L12:M12").Merge()
$ws.Range("N12:O12").Merge()

$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "MEBEFAC 200 MG SR 30 F.C. TABS"
$ws.Range("H12").Value = "0:0"

$fmt = $ws.Range("L12").NumberFormat
$ws.Range("L12").NumberFormat = "@"
$ws.Range("L12").Value = "1"
$ws.Range("L12").NumberFormat = $fmt

$ws.Range("N12").Value = "66.00"

$fmt = $ws.Range("P12").NumberFormat
$ws.Range("P12").NumberFormat = "@"
$ws.Range("P12").Value = "66.0000"
$ws.Range("P12").NumberFormat = $fmt

$ws.Range("Q12").Value = "1:0"

# PRONTOGEST is now at row 13, سرنجات at row 14; renumber PRONTOGEST.
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8

# ------------------------------------------------------------------
# 3) Update the totals row (selling-price sum) - now at row 15.
# ------------------------------------------------------------------
$ws.Range("P15").Value = 332.82999999999998

# ------------------------------------------------------------------
# 4) Update the generated-at timestamp in the footer (now row 16).
# ------------------------------------------------------------------
$ws.Range("A16").Value = "Saturday, 24 May, 2025 9:58 AM"
